$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date line
Replace-Text "2024-01-05 Friday" "2024-01-06 Saturday"

# Row 1
Replace-Text "16×39=624" "72×83=5976"
Replace-Text "57×18=1026" "99×52=5148"
Replace-Text "32×13=416" "35×17=595"
Replace-Text "78×21=1638" "16×72=1152"
Replace-Text "60×49=2940" "20×21=420"

# Row 5
Replace-Text "38×15=570" "59×82=4838"
Replace-Text "93×76=7068" "30×35=1050"
Replace-Text "79×64=5056" "89×29=2581"
Replace-Text "87×15=1305" "28×17=476"
Replace-Text "56×16=896" "87×46=4002"

# Row 10
Replace-Text "60×36=2160" "96×63=6048"
Replace-Text "50×98=4900" "96×46=4416"
Replace-Text "48×17=816" "67×59=3953"
Replace-Text "62×32=1984" "52×78=4056"
Replace-Text "46×93=4278" "43×17=731"

# Row 15 (two cells share "48×89=4272" with row 20 -- disambiguate via table cell)
$t = $d.Tables.Item(1)
Replace-Text "26×49=1274" "23×76=1748"
Replace-Text "78×67=5226" "26×12=312"
Replace-Text "97×46=4462" "53×61=3233"
$t.Cell(15, 4).Range.Text = "99×78=7722"
Replace-Text "98×95=9310" "37×95=3515"

# Row 20
Replace-Text "20×39=780" "13×45=585"
Replace-Text "58×97=5626" "24×29=696"
$t.Cell(20, 3).Range.Text = "39×49=1911"
Replace-Text "73×69=5037" "96×83=7968"
Replace-Text "42×79=3318" "44×23=1012"
